$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Query8 "Affects for Procedures" test to the Basic* naming used elsewhere
$ws.Range("A33").Value = "TestBasicAffectsforProcedures_Query8"

# Insert a new row for the "AffectsStar" test directly below the (renamed) Affects test
$ws.Rows(34).Insert()
$ws.Range("A34").Value = "TestBasicAffectsStarforProcedures_Query8"
$ws.Range("B34").Value = "Source8"

# Rename Query7 basic-modify test to the Modifies-for-statement naming used elsewhere
$ws.Range("A29").Value = "TestBasicModifiesForStmt_Query7"
